# faturamento_anual.xlsx - "correcao no dash e nas modelagens de transacao"
#
# 1) Rename several header labels (B,C,E,F,G,H keep their relative
#    meaning but get the "Faturamento em ..." prefix); the old
#    "Evolução Total (%)" header (I1) is dropped and I1/J1/K1 take over
#    the labels that used to sit one column to the right (Qtd Produtos /
#    Qtd Serviços / Total Itens), while L1/M1/N1 are brand new columns
#    (Qtd Vendas, Ticket Médio Anual, Evolução Ticket Médio (%)).
# 2) The underlying I..L numeric data shifts one column to the right as
#    well (old J->I, old K->J, old L->K) and L gets new "Qtd Vendas"
#    figures; M/N are populated with the new Ticket Médio metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) -------------------------------------------------
$ws.Range("B1").Value = "Faturamento em Produtos"
$ws.Range("C1").Value = "Faturamento em Serviços"
$ws.Range("E1").Value = "Faturamento em Serviços (%)"
$ws.Range("F1").Value = "Faturamento em Produtos (%)"
$ws.Range("G1").Value = "Evolução Faturamento em Serviços (%)"
$ws.Range("H1").Value = "Evolução Faturamento em Produtos (%)"
$ws.Range("I1").Value = "Qtd Produtos"
$ws.Range("J1").Value = "Qtd Serviços"
$ws.Range("K1").Value = "Total Itens"
$ws.Range("L1").Value = "Qtd Vendas"

# M1/N1 are brand-new header cells - give them values first, then pick up
# the same (bold / bordered / centered) header formatting the rest of row 1
# already has, by copying it over from the neighbouring L1 cell.
$ws.Range("M1").Value = "Ticket Médio Anual"
$ws.Range("N1").Value = "Evolução Ticket Médio (%)"

$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("L1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# ---- Row 2 (2018) --------------------------------------------------------
$ws.Range("I2").Value = 1427
$ws.Range("J2").Value = 1322
$ws.Range("K2").Value = 2749
$ws.Range("L2").Value = 2303
$ws.Range("M2").Value = 249.8849326964828

# ---- Row 3 (2019) --------------------------------------------------------
$ws.Range("I3").Value = 17636
$ws.Range("J3").Value = 3034
$ws.Range("K3").Value = 20670
$ws.Range("L3").Value = 13607
$ws.Range("M3").Value = 119.3688689645036
$ws.Range("N3").Value = -52.23046556812919

# ---- Row 4 (2020) --------------------------------------------------------
$ws.Range("I4").Value = 34974
$ws.Range("J4").Value = 3932
$ws.Range("K4").Value = 38906
$ws.Range("L4").Value = 24419
$ws.Range("M4").Value = 116.0415606699701
$ws.Range("N4").Value = -2.787417124244429

# ---- Row 5 (2021) --------------------------------------------------------
$ws.Range("I5").Value = 50119
$ws.Range("J5").Value = 3672
$ws.Range("K5").Value = 53791
$ws.Range("L5").Value = 33150
$ws.Range("M5").Value = 148.2543843137255
$ws.Range("N5").Value = 27.75972975352408

# ---- Row 6 (2022) --------------------------------------------------------
$ws.Range("I6").Value = 57665
$ws.Range("J6").Value = 3867
$ws.Range("K6").Value = 61532
$ws.Range("L6").Value = 38648
$ws.Range("M6").Value = 134.5980547505692
$ws.Range("N6").Value = -9.211417002183008

# ---- Row 7 (2023) --------------------------------------------------------
$ws.Range("I7").Value = 57771
$ws.Range("J7").Value = 3793
$ws.Range("K7").Value = 61564
$ws.Range("L7").Value = 38317
$ws.Range("M7").Value = 118.9145608998617
$ws.Range("N7").Value = -11.65209547773294

# ---- Row 8 (2024) --------------------------------------------------------
$ws.Range("I8").Value = 64032
$ws.Range("J8").Value = 3217
$ws.Range("K8").Value = 67249
$ws.Range("L8").Value = 44870
$ws.Range("M8").Value = 146.472852908402
$ws.Range("N8").Value = 23.174867568781

# ---- Row 9 (2025) --------------------------------------------------------
$ws.Range("I9").Value = 27908
$ws.Range("J9").Value = 1186
$ws.Range("K9").Value = 29094
$ws.Range("L9").Value = 20054
$ws.Range("M9").Value = 162.1813678069213
$ws.Range("N9").Value = 10.724523068
